# Insert a new weekly price record for "Pepino ensalada" (Terminal
# Hortofrutícola Agro Chillán) at row 293, pushing the existing rows
# 293-385 down to 294-386 (dimension grows from A1:R385 to A1:R386).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 293..385 down by one row, creating a blank row 293.
$ws.Rows.Item(293).Insert()

# Populate the new row 293 with the new record.
$ws.Range("A293").Value = 7
$ws.Range("B293").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C293").Value = "Ñuble"
$ws.Range("D293").Value = 45215
$ws.Range("E293").Value = 16
$ws.Range("F293").Value = 100112043
$ws.Range("G293").Value = "Pepino ensalada"
$ws.Range("H293").Value = "Sin especificar"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 100
$ws.Range("K293").Value = 15000
$ws.Range("L293").Value = 15000
$ws.Range("M293").Value = 15000
$ws.Range("N293").Value = "$/caja 60 unidades"
$ws.Range("O293").Value = "Región de Arica y Parinacota"
$ws.Range("P293").Value = 250
$ws.Range("Q293").Value = 60
$ws.Range("R293").Value = "Hortaliza"
